# Update TPM-derived values in Edn1-Ednrb sheet (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.44313933333333
$ws.Range("H2").Value = 37.329418
$ws.Range("I2").Value = 0.9304541596872169
$ws.Range("J2").Value = 0.930454159687217
$ws.Range("M2").Value = 52.271196
$ws.Range("N2").Value = 156.813588
$ws.Range("O2").Value = 0.6500553798777896
$ws.Range("P2").Value = 0.6500553798777895
$ws.Range("Q2").Value = 650.417774947976
$ws.Range("R2").Value = 5853.759974531784
$ws.Range("S2").Value = 0.6048467322343433
$ws.Range("T2").Value = 0.6048467322343433
$ws.Range("G3").Value = 12.44313933333333
$ws.Range("H3").Value = 37.329418
$ws.Range("I3").Value = 0.9304541596872169
$ws.Range("J3").Value = 0.930454159687217
$ws.Range("O3").Value = 0.001694346062422021
$ws.Range("P3").Value = 0.00169434606242202
$ws.Range("Q3").Value = 1.695290632191333
$ws.Range("R3").Value = 15.257615689722
$ws.Range("S3").Value = 0.001576511341730226
$ws.Range("T3").Value = 0.001576511341730226
$ws.Range("G4").Value = 12.44313933333333
$ws.Range("H4").Value = 37.329418
$ws.Range("I4").Value = 0.9304541596872169
$ws.Range("J4").Value = 0.930454159687217
$ws.Range("M4").Value = 24.41792966666667
$ws.Range("N4").Value = 73.253789
$ws.Range("O4").Value = 0.3036664120961408
$ws.Range("P4").Value = 0.3036664120961408
$ws.Range("Q4").Value = 303.8357010738669
$ws.Range("R4").Value = 2734.521309664802
$ws.Range("S4").Value = 0.2825476762921468
$ws.Range("T4").Value = 0.2825476762921468
$ws.Range("G5").Value = 12.44313933333333
$ws.Range("H5").Value = 37.329418
$ws.Range("I5").Value = 0.9304541596872169
$ws.Range("J5").Value = 0.930454159687217
$ws.Range("M5").Value = 3.585005
$ws.Range("N5").Value = 10.755015
$ws.Range("O5").Value = 0.04458386196364773
$ws.Range("P5").Value = 0.04458386196364771
$ws.Range("Q5").Value = 44.60871672569667
$ws.Range("R5").Value = 401.47845053127
$ws.Range("S5").Value = 0.04148323981899672
$ws.Range("T5").Value = 0.04148323981899671
$ws.Range("G6").Value = 0.8272856666666667
$ws.Range("I6").Value = 0.0618615101204856
$ws.Range("J6").Value = 0.06186151012048561
$ws.Range("M6").Value = 52.271196
$ws.Range("N6").Value = 156.813588
$ws.Range("O6").Value = 0.6500553798777896
$ws.Range("P6").Value = 0.6500553798777895
$ws.Range("Q6").Value = 43.24321123032401
$ws.Range("R6").Value = 389.1889010729161
$ws.Range("S6").Value = 0.04021340746118599
$ws.Range("T6").Value = 0.04021340746118599
$ws.Range("G7").Value = 0.8272856666666667
$ws.Range("I7").Value = 0.0618615101204856
$ws.Range("J7").Value = 0.06186151012048561
$ws.Range("O7").Value = 0.001694346062422021
$ws.Range("P7").Value = 0.00169434606242202
$ws.Range("S7").Value = 0.0001048148060881248
$ws.Range("T7").Value = 0.0001048148060881248
$ws.Range("G8").Value = 0.8272856666666667
$ws.Range("I8").Value = 0.0618615101204856
$ws.Range("J8").Value = 0.06186151012048561
$ws.Range("M8").Value = 24.41792966666667
$ws.Range("N8").Value = 73.253789
$ws.Range("O8").Value = 0.3036664120961408
$ws.Range("P8").Value = 0.3036664120961408
$ws.Range("Q8").Value = 20.20060322290811
$ws.Range("R8").Value = 181.805429006173
$ws.Range("S8").Value = 0.01878526282513697
$ws.Range("T8").Value = 0.01878526282513697
$ws.Range("G9").Value = 0.8272856666666667
$ws.Range("I9").Value = 0.0618615101204856
$ws.Range("J9").Value = 0.06186151012048561
$ws.Range("M9").Value = 3.585005
$ws.Range("N9").Value = 10.755015
$ws.Range("O9").Value = 0.04458386196364773
$ws.Range("P9").Value = 0.04458386196364771
$ws.Range("Q9").Value = 2.965823251428334
$ws.Range("R9").Value = 26.692409262855
$ws.Range("S9").Value = 0.002758025028074527
$ws.Range("T9").Value = 0.002758025028074527
$ws.Range("G10").Value = 0.073169
$ws.Range("H10").Value = 0.219507
$ws.Range("I10").Value = 0.005471320266243153
$ws.Range("J10").Value = 0.005471320266243153
$ws.Range("M10").Value = 52.271196
$ws.Range("N10").Value = 156.813588
$ws.Range("O10").Value = 0.6500553798777896
$ws.Range("P10").Value = 0.6500553798777895
$ws.Range("Q10").Value = 3.824631140124
$ws.Range("R10").Value = 34.421680261116
$ws.Range("S10").Value = 0.003556661174105742
$ws.Range("T10").Value = 0.003556661174105742
$ws.Range("G11").Value = 0.073169
$ws.Range("H11").Value = 0.219507
$ws.Range("I11").Value = 0.005471320266243153
$ws.Range("J11").Value = 0.005471320266243153
$ws.Range("O11").Value = 0.001694346062422021
$ws.Range("P11").Value = 0.00169434606242202
$ws.Range("Q11").Value = 0.009968764067000001
$ws.Range("R11").Value = 0.089718876603
$ws.Range("S11").Value = 0.000009270309949358889
$ws.Range("T11").Value = 0.000009270309949358887
$ws.Range("G12").Value = 0.073169
$ws.Range("H12").Value = 0.219507
$ws.Range("I12").Value = 0.005471320266243153
$ws.Range("J12").Value = 0.005471320266243153
$ws.Range("M12").Value = 24.41792966666667
$ws.Range("N12").Value = 73.253789
$ws.Range("O12").Value = 0.3036664120961408
$ws.Range("P12").Value = 0.3036664120961408
$ws.Range("Q12").Value = 1.786635495780333
$ws.Range("R12").Value = 16.079719462023
$ws.Range("S12").Value = 0.00166145619467896
$ws.Range("T12").Value = 0.00166145619467896
$ws.Range("G13").Value = 0.073169
$ws.Range("H13").Value = 0.219507
$ws.Range("I13").Value = 0.005471320266243153
$ws.Range("J13").Value = 0.005471320266243153
$ws.Range("M13").Value = 3.585005
$ws.Range("N13").Value = 10.755015
$ws.Range("O13").Value = 0.04458386196364773
$ws.Range("P13").Value = 0.04458386196364771
$ws.Range("Q13").Value = 0.262311230845
$ws.Range("R13").Value = 2.360801077605
$ws.Range("S13").Value = 0.000243932587509093
$ws.Range("T13").Value = 0.000243932587509093
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.029595
$ws.Range("H14").Value = 0.088785
$ws.Range("I14").Value = 0.002213009926054287
$ws.Range("J14").Value = 0.002213009926054287
$ws.Range("M14").Value = 52.271196
$ws.Range("N14").Value = 156.813588
$ws.Range("O14").Value = 0.6500553798777896
$ws.Range("P14").Value = 0.6500553798777895
$ws.Range("Q14").Value = 1.54696604562
$ws.Range("R14").Value = 13.92269441058
$ws.Range("S14").Value = 0.001438579008154538
$ws.Range("T14").Value = 0.001438579008154538
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.029595
$ws.Range("H15").Value = 0.088785
$ws.Range("I15").Value = 0.002213009926054287
$ws.Range("J15").Value = 0.002213009926054287
$ws.Range("O15").Value = 0.001694346062422021
$ws.Range("P15").Value = 0.00169434606242202
$ws.Range("Q15").Value = 0.004032111585
$ws.Range("R15").Value = 0.036289004265
$ws.Range("S15").Value = 0.000003749604654310929
$ws.Range("T15").Value = 0.000003749604654310927
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.029595
$ws.Range("H16").Value = 0.088785
$ws.Range("I16").Value = 0.002213009926054287
$ws.Range("J16").Value = 0.002213009926054287
$ws.Range("M16").Value = 24.41792966666667
$ws.Range("N16").Value = 73.253789
$ws.Range("O16").Value = 0.3036664120961408
$ws.Range("P16").Value = 0.3036664120961408
$ws.Range("Q16").Value = 0.722648628485
$ws.Range("R16").Value = 6.503837656365
$ws.Range("S16").Value = 0.0006720167841780512
$ws.Range("T16").Value = 0.0006720167841780511
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.029595
$ws.Range("H17").Value = 0.088785
$ws.Range("I17").Value = 0.002213009926054287
$ws.Range("J17").Value = 0.002213009926054287
$ws.Range("M17").Value = 3.585005
$ws.Range("N17").Value = 10.755015
$ws.Range("O17").Value = 0.04458386196364773
$ws.Range("P17").Value = 0.04458386196364771
$ws.Range("Q17").Value = 0.106098222975
$ws.Range("R17").Value = 0.9548840067750001
$ws.Range("S17").Value = 0.00009866452906738658
$ws.Range("T17").Value = 0.00009866452906738656
